$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Housekeeping: drop the stale per-row default formatting on rows 28-29 ---
# (those rows previously carried an unused row-level highlight style that no
# cell actually used any more; clearing it lets the workbook's style table
# collapse the now-orphaned fill/format entry on save). Re-apply the real
# (still in-use) highlight formatting straight back onto the cells so their
# appearance is unchanged.
$ws.Range("A28:K29").Copy() | Out-Null
$ws.Rows.Item(28).ClearFormats()
$ws.Rows.Item(29).ClearFormats()
$ws.Range("A23:K23").Copy() | Out-Null
$ws.Range("A28:K28").PasteSpecial(-4122) | Out-Null
$ws.Range("A29:K29").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- Update fermentation / separation improvement parameters (uncertainty scenarios) ---
# Fermentation acetate loading, midpoint value
$ws.Range("E23").Value = 0.73
# Fermentation CSL loading, midpoint value
$ws.Range("E24").Value = 68
# Fermentation citrate yield, midpoint value
$ws.Range("E28").Value = 0.048

# Re-point the lower/upper bound formulas at the (now updated) midpoint cell
# instead of the previously hard-coded value
$ws.Range("G28").Formula = "=E28*0.0463/0.2087"
$ws.Range("I28").Formula = "=E28*0.34/0.2087"

# --- Restore the author's last active selection ---
$ws.Range("B17").Select()
